# Undo / Redo sequence diagram update ("Update Undo / Redo DG"):
#   - Rename the ":Address" / "BookParser" lifeline label to ":EPiggyParser"
#     (folding its two paragraphs into one, matching the new wording).
#   - Rename "AddressBook" -> "EPiggy" inside the "undoAddressBook()" call-out.
#   - Rename "ReadOnlyAddressBook" -> "ReadOnlyEPiggy" inside the
#     "resetData(ReadOnlyAddressBook)" call-out.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate shapes by their current text rather than hard-coded indices, so the
# script is resilient to shape-ordering differences.
$lifeline = $null
$undoShape = $null
$resetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }
    $t = $shp.TextFrame.TextRange.Text
    if ($t.Contains("BookParser")) { $lifeline = $shp }
    elseif ($t.Contains("undoAddressBook")) { $undoShape = $shp }
    elseif ($t.Contains("ReadOnlyAddressBook")) { $resetShape = $shp }
}

# --- 1) ":Address" / "BookParser" (two paragraphs) -> ":" + "EPiggyParser" (one paragraph) ---
$lifelineRange = $lifeline.TextFrame.TextRange

# Re-writing the whole text range merges the two paragraphs into a single
# paragraph ":EPiggyParser" (the surviving paragraph mark/properties come
# from the former second paragraph).
$lifelineRange.Text = ":EPiggyParser"

# Re-apply run-level formatting to the "EPiggyParser" portion (chars 2-13) so
# it stays a distinct run, separate from the leading ":" run. Both runs keep
# the shape's original size/colour; touching Font.Size is enough to force the
# run split without disturbing the inherited (scheme-coloured) solid fill.
$parserRun = $lifelineRange.Characters(2, 12)
$parserRun.Font.Size = 16

# --- 2) "undoAddressBook()" -> "undoEPiggy()" ---
$undoRange = $undoShape.TextFrame.TextRange
$undoText = $undoRange.Text
$idx = $undoText.IndexOf("AddressBook")
$undoRange.Characters($idx + 1, "AddressBook".Length).Text = "EPiggy"

# --- 3) "resetData(ReadOnlyAddressBook)" -> "resetData(ReadOnlyEPiggy)" ---
$resetRange = $resetShape.TextFrame.TextRange
$resetText = $resetRange.Text
$idx2 = $resetText.IndexOf("ReadOnlyAddressBook")
$resetRange.Characters($idx2 + 1, "ReadOnlyAddressBook".Length).Text = "ReadOnlyEPiggy"
